# Blackjack_Specification.pptx — "Add TestSheet Change Source.c"
#
# 1) Slide 1 subtitle: "Required specifications" -> "specifications"
#    (author selected/removed the leading "Required " word).
# 2) The cached text of every "datetimeFigureOut" date field on the
#    slide master and all custom layouts advances from 2019/3/17 to
#    2019/3/29 (PowerPoint re-stamps these automatic date placeholders
#    whenever the deck is re-saved on a later day).

$p = $ppt.ActivePresentation

# --- 1) Slide 1 subtitle text -------------------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shape = $slide1.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -eq "Required specifications") {
            $tr = $shape.TextFrame.TextRange
            $tr.Delete()
            $tr.InsertAfter("specifications")
        }
    }
}

# --- 2) Refresh the cached "datetimeFigureOut" placeholders -------------
function Update-DatePlaceholder {
    param($shapes, $newText)

    $updated = $false
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        $phType = $null
        try { $phType = $shape.PlaceholderFormat.Type } catch {}
        if ($phType -eq 16) {
            if ($shape.TextFrame.TextRange.Text -ne $newText) {
                $shape.TextFrame.TextRange.Text = $newText
            }
            $updated = $true
        }
    }
    return $updated
}

$newDate = "2019/3/29"
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes $newDate | Out-Null

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes $newDate | Out-Null
}
